# Weekly fruit/vegetable price update: insert a new week's record for
# "Coliflor" at Terminal Hortofrutícola Agro Chillán, pushing the existing
# rows 424..535 down to 425..536.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 424; Excel shifts rows 424:535 down to
# 425:536 and grows the used range to A1:R536 automatically.
$ws.Rows(424).Insert()

# Populate the newly inserted row 424 with the new weekly record.
$ws.Cells.Item(424, 1).Value  = 7
$ws.Cells.Item(424, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(424, 3).Value  = "Ñuble"
$ws.Cells.Item(424, 4).Value  = 45135
$ws.Cells.Item(424, 5).Value  = 16
$ws.Cells.Item(424, 6).Value  = 100112008
$ws.Cells.Item(424, 7).Value  = "Coliflor"
$ws.Cells.Item(424, 8).Value  = "Sin especificar"
$ws.Cells.Item(424, 9).Value  = "Primera"
$ws.Cells.Item(424, 10).Value = 180
$ws.Cells.Item(424, 11).Value = 1000
$ws.Cells.Item(424, 12).Value = 1000
$ws.Cells.Item(424, 13).Value = 1000
$ws.Cells.Item(424, 14).Value = "$/unidad"
$ws.Cells.Item(424, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(424, 16).Value = 1000
$ws.Cells.Item(424, 17).Value = 1
$ws.Cells.Item(424, 18).Value = "Hortaliza"
